# Auto-generated edit script applying the diff changes
$wb = $excel.ActiveWorkbook

# sheet1 (展览)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2513
$ws1.Range("F3").Value = 541
$ws1.Range("F4").Value = 449
$ws1.Range("F5").Value = 279
$ws1.Range("F10").Value = 282
$ws1.Range("F13").Value = 5413
$ws1.Range("F15").Value = 1611
$ws1.Range("F16").Value = 3957
$ws1.Range("F17").Value = 392
$ws1.Range("F20").Value = 4460
$ws1.Range("F21").Value = 5859
$ws1.Range("F24").Value = 642
$ws1.Range("F25").Value = 3634
$ws1.Range("F29").Value = 115
$ws1.Range("F30").Value = 951
$ws1.Range("F31").Value = 1334
$ws1.Range("F32").Value = 396
$ws1.Range("F33").Value = 448
$ws1.Range("F34").Value = 1540
$ws1.Range("F36").Value = 1612
$ws1.Range("F37").Value = 146
$ws1.Range("F38").Value = 1058
$ws1.Range("F41").Value = 597
$ws1.Range("F43").Value = 182
$ws1.Range("F44").Value = 2721
$ws1.Range("F49").Value = 3842

# sheet2 (演出)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 1148
$ws2.Range("F14").Value = 10

# sheet3 (本地生活)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 3612

# sheet4 (全部类型)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("B2").Value = "'2024-06-15"
$ws4.Range("C2").Value = "杭州·《菊次郎的夏天》久石让宫崎骏经典作品主题音乐会"
$ws4.Range("D2").Value = "武林路77号 浙江省文化馆小剧场（原群艺馆小剧场）"
$ws4.Range("E2").Value = "2024.06.15 19:45-06.15 21:00"
$ws4.Range("F2").Value = 4
$ws4.Range("G2").Value = 100
$ws4.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=86753"
$ws4.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202405/4PeKYE9D1717136762813.jpeg"
$ws4.Range("F3").Value = 2513
$ws4.Range("F4").Value = 541
$ws4.Range("F5").Value = 449
$ws4.Range("F6").Value = 279
$ws4.Range("F7").Value = 1148
$ws4.Range("F12").Value = 282
$ws4.Range("F15").Value = 5414
$ws4.Range("F17").Value = 1611
$ws4.Range("F18").Value = 4460
$ws4.Range("F19").Value = 5859
$ws4.Range("F22").Value = 642
$ws4.Range("F23").Value = 3634
$ws4.Range("F27").Value = 115
$ws4.Range("C28").Value = "杭州·夏之誓国乙only-日夜场"
$ws4.Range("D28").Value = "北干街道萧杭路689号 杭州时尚外滩艺术中心"
$ws4.Range("E28").Value = "2024.07.27 10:00-07.27 21:00"
$ws4.Range("F28").Value = 1334
$ws4.Range("G28").Value = 99
$ws4.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=83589"
$ws4.Range("I28").Value = "//i2.hdslb.com/bfs/openplatform/202405/99kWb2dy1714964533903.png"
$ws4.Range("C29").Value = "杭州·文豪野犬only"
$ws4.Range("D29").Value = "望江东路333号 杭州瑞莱克斯大酒店"
$ws4.Range("E29").Value = "2024.07.27 10:00-07.27 17:00"
$ws4.Range("F29").Value = 397
$ws4.Range("G29").Value = 60
$ws4.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=86859"
$ws4.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202406/LQOrzPac1717473481789.png"
$ws4.Range("C30").Value = "杭州·第五人格only"
$ws4.Range("F30").Value = 449
$ws4.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=86861"
$ws4.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202406/ZiqYUjIt1717471158198.jpeg"
$ws4.Range("C31").Value = "杭州·黄西全新脱口秀专场《水土不服》"
$ws4.Range("D31").Value = "延安路279号 浙江胜利剧院"
$ws4.Range("E31").Value = "2024.07.27 19:30-07.27 21:30"
$ws4.Range("F31").Value = 2
$ws4.Range("G31").Value = 224
$ws4.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=85935"
$ws4.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202405/9YqhR4Ke1716191781870.jpeg"
$ws4.Range("F32").Value = 1540
$ws4.Range("F34").Value = 1612
$ws4.Range("F36").Value = 1058
$ws4.Range("F38").Value = 597
$ws4.Range("F43").Value = 2721
$ws4.Range("F49").Value = 3843
